$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 already holds "システム設計の面接試験" data; use it as the format
# template for the two new rows being appended below it (so the new rows
# pick up the sheet's standard cell style instead of a blank default one).
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A11:E11").PasteSpecial(-4122)

# New row 10: another Kindle book purchase, still undecided.
$ws.Range("A10").Value = "ハンズオンで学ぶAWSコスト最適化入門"
$ws.Range("B10").Value = "緒方遼太郎"
$ws.Range("C10").Value = 2000
$ws.Range("D10").Value = "Kindle"
$ws.Range("E10").Value = "未定"

# New row 11: a test book entry.
$ws.Range("A11").Value = "テスト書籍"
$ws.Range("B11").Value = "テスト著者"
$ws.Range("C11").Value = 1000
$ws.Range("D11").Value = "単行本"
$ws.Range("E11").Value = "未定"
